$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: best_params updates
$ws.Range("B2").Value = "{'alpha': 0.1, 'max_iter': 1000}"
$ws.Range("F2").Value = "{'max_depth': 50, 'n_estimators': 200}"
$ws.Range("G2").Value = "{'learning_rate': 0.1, 'max_depth': 5, 'n_estimators': 100}"
$ws.Range("H2").Value = "{'learning_rate': 0.1, 'n_estimators': 200}"
$ws.Range("J2").Value = "{'learning_rate': 0.1, 'max_depth': 5, 'n_estimators': 100}"
$ws.Range("K2").Value = "{'activation': 'leaky_relu', 'b_random_vec_range': [0, 10], 'lam': 2, 'n_layer': 64, 'n_nodes': 128, 'random_seed': 882, 'same_feature': True, 'w_random_vec_range': [-10, 10]}"

# Row 3: rmse
$ws.Range("B3").Value = 0.0885986085186389
$ws.Range("C3").Value = 0.1104093246962435
$ws.Range("D3").Value = 0.09813628659031297
$ws.Range("E3").Value = 0.0978653792501606
$ws.Range("F3").Value = 0.07805456951742873
$ws.Range("G3").Value = 0.07692894996935146
$ws.Range("H3").Value = 0.1106237932078655
$ws.Range("I3").Value = 0.0738026804302649
$ws.Range("J3").Value = 0.07278482903015036
$ws.Range("K3").Value = 0.04492345241591222

# Row 4: r2
$ws.Range("B4").Value = 0.7231896947312604
$ws.Range("C4").Value = 0.5857057056389258
$ws.Range("D4").Value = 0.669936750379415
$ws.Range("E4").Value = 0.6367001382937567
$ws.Range("F4").Value = 0.791662064613185
$ws.Range("G4").Value = 0.7976680914211356
$ws.Range("H4").Value = 0.5429244818562209
$ws.Range("I4").Value = 0.7742791726891379
$ws.Range("J4").Value = 0.8195502283386082
$ws.Range("K4").Value = 0.921595710794913

# Row 5: mape
$ws.Range("B5").Value = 27.47949770088503
$ws.Range("C5").Value = 56.30056283062441
$ws.Range("D5").Value = 36.22051711876593
$ws.Range("E5").Value = 62.71131990337974
$ws.Range("F5").Value = 20.40586830188068
$ws.Range("G5").Value = 21.83669234854217
$ws.Range("H5").Value = 84.6733506131568
$ws.Range("I5").Value = 19.4730159888745
$ws.Range("J5").Value = 25.26368650063011
$ws.Range("K5").Value = 10.56005703721227
